$d = $word.ActiveDocument

# Delete everything after the Title paragraph (paragraph 1), keeping the
# sectPr (attached to the body) intact.
$titleEnd = $d.Paragraphs(1).Range.End
$docEnd = $d.Content.End
$r = $d.Range($titleEnd, $docEnd)
$r.Delete()

$insertionPoint = $d.Range($titleEnd, $titleEnd)

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t>Project Members: Cam Foster, Anis Ali, Michael Alread, Ida Astaneh</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t>Release Date: 2020-10-27</w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>The Original Data</w:t></w:r></w:p>
<w:p><w:r><w:t>&lt;Need to add stuff here&gt;</w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Process Overview</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">To expand our original dataset and get the audio features for each track we needed clean the dataset to filter for only songs in the US region for 2017 extract just the unique track IDs. That allowed us to call the API fewer times and reduce the time it took to get our data. </w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Once we had the audio features for all songs in our final </w:t></w:r><w:r><w:t>dataset,</w:t></w:r><w:r><w:t xml:space="preserve"> we merged the original dataset with the audio features dataset on track ID</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">giving us our final, working dataset. </w:t></w:r></w:p>
<w:p><w:r><w:t>&lt;Talk about the analysis process&gt;</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Cleaning the Data</w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Spotify API</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">There is a python library called spotipy that is used to interact with the Spotify Web API. </w:t></w:r><w:r><w:t xml:space="preserve">You are required to create a developer account and get a Client ID and Client Secret key to authenticate with Spotify. Most of the functions in this library are geared towards interacting with Spotify for web applications and doing things such as selecting the next track, getting related artists, and get current user data. We used two functions, </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>spotipy.track</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">() and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>spotipy.audio_features</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(), to get information about the track requested. </w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">A function was written to call the API and return the track identifiers (i.e. artist, track name, album name, etc.) and audio features (i.e. loudness, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>acousticness</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, energy, etc.). </w:t></w:r><w:r><w:t>We ran into an issue with the API timing out but adding a pause under an except clause gave us a work around for this issue.</w:t></w:r></w:p>
<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>All of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the audio features are defined by Spotify on their developer website. </w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Questions and Hypothesis</w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t xml:space="preserve">Given the data provided, is there a direct correlation between our Spotipy headers? </w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Could we bucket the stream counts to show the percentage breakdown? How &#8220;popular&#8221; is popular? </w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t xml:space="preserve">Which artist/album appears the most and least amount of times on the list? </w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t xml:space="preserve">What decade, prior to 2010, was responsible for providing the most streams in 2017? </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertionPoint.InsertXML($xml)
